$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")

# Append new row 78: Raw and Clean from SSA for July 17th (2020-08-16 pull)
# Force column A to be stored as text (like the rest of the date column)
# instead of letting Excel auto-convert the ISO date string into a date
# serial number.
$ws.Range("A78").NumberFormat = "@"
$ws.Range("A78").Value = "2020-08-16"
$ws.Range("A78").Style = "Normal"

$ws.Range("B78").Value = 522162
$ws.Range("C78").Value = 573723
$ws.Range("D78").Value = 81046
$ws.Range("E78").Value = 56757
$ws.Range("F78").Value = 26.21
